$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 40, shifting the existing rows
# 40-47 down to 41-48 (new weekly price observation for Espárragos).
$ws.Rows("40:40").Insert()

$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44841
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = 300000000
$ws.Cells.Item(40, 7).Value = "Espárragos"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 300
$ws.Cells.Item(40, 11).Value = 1800
$ws.Cells.Item(40, 12).Value = 1800
$ws.Cells.Item(40, 13).Value = 1800
$ws.Cells.Item(40, 14).Value = "$/kilo"
$ws.Cells.Item(40, 15).Value = "Provincia de Linares"
$ws.Cells.Item(40, 16).Value = 1800
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
